$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.3
$ws.Range("E3").Value = 1.27
$ws.Range("B4").Value = 1.42
$ws.Range("F4").Value = 1.08
$ws.Range("C5").Value = 1.37
$ws.Range("D5").Value = 1.31
$ws.Range("G5").Value = 0.71
$ws.Range("D6").Value = 1.53
$ws.Range("E7").Value = 1.92
$ws.Range("F7").Value = 1.51
